$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 267, shifting the existing row 267 (and everything below) down to 268.
$ws.Rows.Item(267).Insert()

# Populate the newly inserted row 267 with the new weekly price entry.
$ws.Cells.Item(267, 1).Value = 3
$ws.Cells.Item(267, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(267, 3).Value = "Coquimbo"
$ws.Cells.Item(267, 4).Value = 45218
$ws.Cells.Item(267, 4).NumberFormat = $ws.Cells.Item(268, 4).NumberFormat
$ws.Cells.Item(267, 5).Value = 5
$ws.Cells.Item(267, 6).Value = 100112026
$ws.Cells.Item(267, 7).Value = "Haba"
$ws.Cells.Item(267, 8).Value = "Sin especificar"
$ws.Cells.Item(267, 9).Value = "Primera"
$ws.Cells.Item(267, 10).Value = 65
$ws.Cells.Item(267, 11).Value = 12000
$ws.Cells.Item(267, 12).Value = 12000
$ws.Cells.Item(267, 13).Value = 12000
$ws.Cells.Item(267, 14).Value = '$/saco 25 kilos'
$ws.Cells.Item(267, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(267, 16).Value = 480
$ws.Cells.Item(267, 17).Value = 25
$ws.Cells.Item(267, 18).Value = "Hortaliza"

Write-Host "Row inserted and populated"
